$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 13 values (mean_temp, max_temp, min_temp)
$ws.Range("C13").Value = -5.23035714285714
$ws.Range("D13").Value = 6.8
$ws.Range("E13").Value = -12.2

# Clear row 17 values for mean_temp, max_temp, min_temp
$ws.Range("C17:E17").ClearContents()
